$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.973.55'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.846.02'
$ws.Range("E4").Value = '  +0.11%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.79'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("E7").Value = '  +0.08%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.92'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  +7.25%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.329'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +3.15%  '
$ws.Range("E10").Value = '  +2.20%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("E12").Value = '  +2.12%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.45'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +5.40%  '
$ws.Range("D14").Value = '1.849.97'
$ws.Range("E14").Value = '  +2.38%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.68'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  +3.05%  '
$ws.Range("D17").Value = '34.967.09'
$ws.Range("E17").Value = '  +0.33%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.05'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '0.0₃0793'
$ws.Range("E19").Value = '  +1.49%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.77'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.65%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.23'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +4.13%  '
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  +1.04%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.55'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -0.60%  '
$ws.Range("E26").Value = '  +0.82%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.52'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +1.96%  '
$ws.Range("E28").Value = '  +3.98%  '
$ws.Range("E29").Value = '  +8.94%  '
$ws.Range("E30").Value = '  +0.02%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0554'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +1.47%  '
$ws.Range("E32").Value = '  -0.58%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.94'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  +21.91%  '
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.96'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +11.78%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.747'
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +9.47%  '
$ws.Range("E37").Value = '  +5.32%  '
$ws.Range("E38").Value = '  +11.58%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '89.77'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").Value = '1.349.92'
$ws.Range("E40").Value = '  +3.06%  '
$ws.Range("E41").Value = '  +2.81%  '
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.59'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("E43").Value = '  +4.64%  '
$ws.Range("E44").Value = '  -1.68%  '
$ws.Range("E45").Value = '  +2.41%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0531'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +4.15%  '
$ws.Range("E47").Value = '  +3.93%  '
$ws.Range("E48").Value = '  +1.82%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").Value = '  +15.80%  '
$ws.Range("E51").Value = '  -0.43%  '
